$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A6 previously held "${item.nameAndDisplayValues}" - update the text so it now
# reads the localized label + display values helper.
$ws.Range("A6").Value = '${item.getLocalizedLabelAndDisplayValues(locale)}'

# B5 keeps the same displayed value "${cityId.value}" (it is re-used, this is
# effectively a no-op on content but we set it explicitly for safety).
$ws.Range("B5").Value = '${cityId.value}'

# Update the selected/active cell shown in the sheet view to A7.
$ws.Range("A7").Select()

$wb.Save()
